# "Some more data collection" - remove the scratch UNIQUE() helper column (M),
# hide duplicate rows that are no longer needed for visible reporting, drop
# the unused trailing placeholder rows (35-39), and restore the view to the
# top of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the ad-hoc dynamic-array helper column (M1:M22 held
# =UNIQUE(G:G) plus its cached spill values).
$ws.Columns.Item(13).Delete()

# Remove the trailing placeholder rows (A35:A40 held bare counters 34-39)
# that aren't part of the collected data anymore.
$ws.Range("A35:A40").EntireRow.Delete()

# Hide the rows whose data is a repeat of an already-represented
# configuration, keeping only the distinct / highlighted rows visible.
$hiddenRows = @(2, 4, 5, 7, 8, 11, 13, 14, 15, 16, 17, 19, 20, 21, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# Reset the frozen-pane scroll position back to the top of the data and
# update the saved selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
[void]$ws.Range("G46").Select()
